# ultima modificacion hecha a repo existente
#
# The document ends with two empty paragraphs right before the section
# break: one carrying w14:paraId "7F137BED" and the final one carrying
# w14:paraId "50B20278". This inserts a brand-new paragraph containing
# the text "ultima modifiocacion final" between those two, matching the
# 7F137BED paragraph's formatting (which InsertParagraphAfter naturally
# inherits).

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$lastEmpty = $d.Paragraphs.Item($count - 1)

$lastEmpty.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($count)
$newPara.Range.Text = "ultima modifiocacion final"
